$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“핸즈온 머신러닝 2판”, “Do It! 딥러닝 입문” 주피터 노트북 업데이트 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/02/18/%ed%95%b8%ec%a6%88%ec%98%a8-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-2%ed%8c%90-do-it-%eb%94%a5%eb%9f%ac%eb%8b%9d-%ec%9e%85%eb%ac%b8-%ec%a3%bc%ed%94%bc%ed%84%b0-%eb%85%b8%ed%8a%b8%eb%b6%81/"

$ws.Range("D23").Value = "안녕하세요, Cognex Research Engineer 이호성입니다. `n요즘 컴퓨터 비전계를 뜨겁게 달구고 있는 모델이 있습니다. 바로 자연"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2719"

$ws.Range("D28").Value = "NVIDIA Jeston 환경 셋팅 2-2 편 (OpenCV 4.1.0 설치 on Jeston Nano)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/88"

$ws.Range("D30").Value = "2020년"

$ws.Range("D37").Value = "[Paper Review] FreeLB: Enhanved Adversarial Training for Natural Language Understanding"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1441&mod=document&pageid=1"

$ws.Range("D39").Value = "Artistic Style Transfer with Deep Learning"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Artistic-Style-Transfer-with-Deep-Learning-1"

$ws.Range("D51").Value = "[독후감] 내가 오늘 하고 있는 활동들을 영어로 표현하지 못한다면, `"야나두 나의 즐거운 영어생활`""
$ws.Range("E51").Value = "https://bskyvision.com/1114"
